$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.471.65'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.59%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.105.89'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.17%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.006'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.61%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '334.15'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.79%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.005'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.59%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5221'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.03%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4536'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +5.09%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '53.45'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +14.25%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.08927'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.06%  '

$ws.Range("E11").Value = '  +1.84%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '24.11'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.38%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '2.101.55'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.92%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.817'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.18%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '8.032'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +4.64%  '

$ws.Range("E16").Value = '  +1.39%  '

$ws.Range("E17").Value = '  +1.58%  '

$ws.Range("E18").Value = '  +0.65%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06656'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.40%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '19.25'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.90%  '

$ws.Range("E21").Value = '  +0.54%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.338'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.42%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '30.529.33'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.56%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.47'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.97%  '

$ws.Range("E25").Value = '  +2.40%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.349.11'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.91%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '22.25'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.51%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '162.84'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.60%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.523'
$ws.Range("D29").Style = "Normal"

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '133.46'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.54%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.208'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.20%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.1072'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.14%  '

$ws.Range("E33").Value = '  +0.14%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.412'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +3.68%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.943'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.97%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '10.42'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +4.61%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.811'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +6.71%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02584'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.62%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06840'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +2.68%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.2295'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.37%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '12.72'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.28%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.6873'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.57%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.246'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.03%  '

$ws.Range("B44").Value = 'NEARProtocol'
$ws.Range("C44").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.313'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +4.89%  '

$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '14.02'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.60%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.6363'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.37%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.661'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.59%  '

$ws.Range("E48").Value = '  +21.57%  '

$ws.Range("E49").Value = '  -0.15%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '83.42'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.31%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.205'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.33%  '
